$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G10").Value = 2.45
$ws.Range("AA10").Value = 26
$ws.Range("AL10").Value = 51
$ws.Range("AP10").Value = 34
$ws.Range("AR10").Value = 101
$ws.Range("AX10").Value = 21
$ws.Range("BA10").Value = 126
$ws.Range("G13").Value = 1.65
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 5.1
$ws.Range("J13").Value = 2.15
$ws.Range("K13").Value = 2.15
$ws.Range("L13").Value = 5.2
$ws.Range("M13").Value = 1.01
$ws.Range("N13").Value = 7.8
$ws.Range("O13").Value = 1.31
$ws.Range("P13").Value = 2.87
$ws.Range("Q13").Value = 1.93
$ws.Range("R13").Value = 1.7
$ws.Range("U13").Value = 1.9
$ws.Range("V13").Value = 1.72
$ws.Range("W13").Value = 6.1
$ws.Range("X13").Value = 7.1
$ws.Range("Y13").Value = 8.25
$ws.Range("Z13").Value = 12
$ws.Range("AA13").Value = 14
$ws.Range("AB13").Value = 30
$ws.Range("AC13").Value = 9
$ws.Range("AD13").Value = 6.9
$ws.Range("AE13").Value = 17.5
$ws.Range("AF13").Value = 90
$ws.Range("AG13").Value = 12.5
$ws.Range("AH13").Value = 30
$ws.Range("AI13").Value = 16.5
$ws.Range("AJ13").Value = 100
$ws.Range("AK13").Value = 55
$ws.Range("AM13").Value = 800
$ws.Range("AN13").Value = 3.4
$ws.Range("AO13").Value = 7.7
$ws.Range("AP13").Value = 17
$ws.Range("AQ13").Value = 25
$ws.Range("AR13").Value = 55
$ws.Range("AS13").Value = 250
$ws.Range("AT13").Value = 2.57
$ws.Range("AU13").Value = 7.5
$ws.Range("AV13").Value = 70
$ws.Range("AW13").Value = 6.7
$ws.Range("AX13").Value = 29
$ws.Range("AY13").Value = 35
$ws.Range("BA13").Value = 200
$ws.Range("BB13").Value = 450
